$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Work_Mode",
    "Job_Location",
    "LinkedIn_Poster",
    "LinkedIn_Posted",
    "Resume_received",
    "Resume_downloaded"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# Build the header format on the first cell, then fan it out to the rest of
# the row with a single format-only paste so every header cell converges on
# the same style index (instead of each cell racking up its own incremental
# style combinations).
$seed = $ws.Cells.Item(1, 1)
$seed.Font.Bold = $true
$seed.HorizontalAlignment = -4108  # xlCenter
$seed.VerticalAlignment = -4160    # xlTop
$seed.Borders.LineStyle = 1        # xlContinuous

$seed.Copy()
$ws.Range("B1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data row
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"
$ws.Range("C2").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Hybrid"
$ws.Range("G2").Value = "Bengaluru, Karnataka, India"

# The multi-line job description otherwise triggers an auto row-height
# calculation; re-running AutoFit restores the default (non-custom) height.
$ws.Rows.Item(2).EntireRow.AutoFit()
